$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 276 (shifts the existing rows 276:342 down to 277:343,
# carrying their values/styles with them automatically).
$ws.Rows.Item(276).Insert()

# Populate the newly inserted row 276 with the new observation.
# (Same market/region/product/category/variety/unit as the rest of the
#  Frambuesa block; only Fecha, Calidad, Volumen differ.)
$ws.Cells.Item(276, 1).Value  = 6
$ws.Cells.Item(276, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(276, 3).Value  = "Metropolitana"
$ws.Cells.Item(276, 4).Value  = 45015
$ws.Cells.Item(276, 5).Value  = 13
$ws.Cells.Item(276, 6).Value  = "Fruta"
$ws.Cells.Item(276, 7).Value  = 100101
$ws.Cells.Item(276, 8).Value  = "Berries"
$ws.Cells.Item(276, 9).Value  = 100101004
$ws.Cells.Item(276, 10).Value = "Frambuesa"
$ws.Cells.Item(276, 11).Value = "Sin especificar"
$ws.Cells.Item(276, 12).Value = "Primera"
$ws.Cells.Item(276, 13).Value = 250
$ws.Cells.Item(276, 14).Value = 7000
$ws.Cells.Item(276, 15).Value = 7000
$ws.Cells.Item(276, 16).Value = 7000
$ws.Cells.Item(276, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(276, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(276, 19).Value = 3500
$ws.Cells.Item(276, 20).Value = 2
